# Apply cryptos list update (price & volume columns) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preventing Excel from auto-converting
# numeric-looking strings (e.g. "1.00", "0.170") into actual numbers, which
# would silently drop formatting such as trailing zeros.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "63.384.79"
$ws.Range("E2").Value = "  +0.11%  "
Set-TextValue "D3" "2.561.38"
$ws.Range("E3").Value = "  +5.22%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.17%  "
Set-TextValue "D5" "569.41"
$ws.Range("E5").Value = "  +0.59%  "
Set-TextValue "D6" "148.54"
$ws.Range("E6").Value = "  +5.04%  "
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.17%  "
Set-TextValue "D8" "0.583"
$ws.Range("E8").Value = "  -1.38%  "
Set-TextValue "D9" "2.565.87"
$ws.Range("E9").Value = "  +5.29%  "
$ws.Range("E10").Value = "  +0.17%  "
Set-TextValue "D11" "5.65"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("E13").Value = "  +0.70%  "
Set-TextValue "D14" "27.41"
$ws.Range("E14").Value = "  +4.35%  "
Set-TextValue "D15" "3.017.59"
$ws.Range("E15").Value = "  +5.06%  "
Set-TextValue "D16" "63.261.47"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  -0.35%  "
Set-TextValue "D18" "2.556.24"
$ws.Range("E18").Value = "  +5.23%  "
Set-TextValue "D19" "11.57"
$ws.Range("E19").Value = "  +2.79%  "
Set-TextValue "D20" "336.69"
$ws.Range("E20").Value = "  -1.25%  "
Set-TextValue "D21" "4.32"
$ws.Range("E21").Value = "  +1.51%  "
Set-TextValue "D22" "6.79"
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("E23").Value = "  +0.30%  "
Set-TextValue "D24" "65.51"
$ws.Range("E24").Value = "  +0.28%  "
Set-TextValue "D25" "0.170"
$ws.Range("E25").Value = "  -2.74%  "
Set-TextValue "D26" "1.60"
$ws.Range("E26").Value = "  +3.97%  "
Set-TextValue "D27" "1.49"
$ws.Range("E27").Value = "  +10.92%  "
$ws.Range("E28").Value = "  -0.13%  "
Set-TextValue "D29" "8.43"
$ws.Range("E29").Value = "  +2.88%  "
Set-TextValue "D30" "7.21"
$ws.Range("E30").Value = "  +8.63%  "
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("E32").Value = "  +1.89%  "
Set-TextValue "D33" "177.54"
$ws.Range("E33").Value = "  +2.01%  "
Set-TextValue "D34" "1.60"
$ws.Range("E34").Value = "  +7.24%  "
Set-TextValue "D35" "415.75"
$ws.Range("E35").Value = "  +12.02%  "
Set-TextValue "D36" "0.399"
$ws.Range("E36").Value = "  -0.08%  "
Set-TextValue "D37" "18.91"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("E40").Value = "  +4.24%  "
Set-TextValue "D41" "1.00"
$ws.Range("E41").Value = "  +0.19%  "
Set-TextValue "D42" "39.43"
$ws.Range("E42").Value = "  -1.28%  "
Set-TextValue "D43" "152.44"
$ws.Range("E43").Value = "  +2.71%  "
Set-TextValue "D44" "3.75"
$ws.Range("E44").Value = "  +1.32%  "
Set-TextValue "D45" "20.72"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("E46").Value = "  +2.38%  "
Set-TextValue "D47" "0.0969"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("E49").Value = "  +5.00%  "
Set-TextValue "D50" "18.49"
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("E51").Value = "  +3.44%  "
